$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("J2").Value = 0.7825479339666588
$ws.Range("M2").Value = 4.618552666666667
$ws.Range("N2").Value = 13.855658
$ws.Range("O2").Value = 0.78434648953826
$ws.Range("P2").Value = 0.78434648953826
$ws.Range("Q2").Value = 0.519822721186
$ws.Range("R2").Value = 4.678404490674
$ws.Range("S2").Value = 0.6137887249021671
$ws.Range("T2").Value = 0.613788724902167

# Row 3 updates
$ws.Range("J3").Value = 0.7825479339666588
$ws.Range("O3").Value = 0.1153620112191035
$ws.Range("P3").Value = 0.1153620112191036
$ws.Range("R3").Value = 0.6881016980879999
$ws.Range("S3").Value = 0.09027630353774799
$ws.Range("T3").Value = 0.090276303537748

# Row 4 updates
$ws.Range("J4").Value = 0.7825479339666588
$ws.Range("M4").Value = 0.478937
$ws.Range("N4").Value = 1.436811
$ws.Range("O4").Value = 0.0813355572127976
$ws.Range("P4").Value = 0.08133555721279762
$ws.Range("Q4").Value = 0.053904838287
$ws.Range("R4").Value = 0.485143544583
$ws.Range("S4").Value = 0.06364897225490175
$ws.Range("T4").Value = 0.06364897225490175

# Row 5 updates
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 0.112551
$ws.Range("H5").Value = 0.337653
$ws.Range("I5").Value = 0.7825479339666589
$ws.Range("J5").Value = 0.7825479339666588
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1116203333333333
$ws.Range("N5").Value = 0.334861
$ws.Range("O5").Value = 0.01895594202983873
$ws.Range("P5").Value = 0.01895594202983874
$ws.Range("Q5").Value = 0.012562980137
$ws.Range("R5").Value = 0.113066821233
$ws.Range("S5").Value = 0.01483393327184205
$ws.Range("T5").Value = 0.01483393327184205

# Row 6 updates
$ws.Range("D6").Value = "ECs"
$ws.Range("G6").Value = 0.03127533333333334
$ws.Range("H6").Value = 0.09382600000000001
$ws.Range("I6").Value = 0.2174520660333412
$ws.Range("J6").Value = 0.2174520660333412
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.618552666666667
$ws.Range("N6").Value = 13.855658
$ws.Range("O6").Value = 0.78434648953826
$ws.Range("P6").Value = 0.78434648953826
$ws.Range("Q6").Value = 0.1444467741675556
$ws.Range("R6").Value = 1.300020967508
$ws.Range("S6").Value = 0.170557764636093
$ws.Range("T6").Value = 0.170557764636093

# Row 7 updates
$ws.Range("D7").Value = "FAPs"
$ws.Range("G7").Value = 0.03127533333333334
$ws.Range("H7").Value = 0.09382600000000001
$ws.Range("I7").Value = 0.2174520660333412
$ws.Range("J7").Value = 0.2174520660333412
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.6792986666666666
$ws.Range("N7").Value = 2.037896
$ws.Range("O7").Value = 0.1153620112191035
$ws.Range("P7").Value = 0.1153620112191036
$ws.Range("Q7").Value = 0.02124529223288889
$ws.Range("R7").Value = 0.191207630096
$ws.Range("S7").Value = 0.02508570768135554
$ws.Range("T7").Value = 0.02508570768135555

# Row 8 (new)
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Fgf5"
$ws.Range("C8").Value = "Fgfr3"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.03127533333333334
$ws.Range("H8").Value = 0.09382600000000001
$ws.Range("I8").Value = 0.2174520660333412
$ws.Range("J8").Value = 0.2174520660333412
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.478937
$ws.Range("N8").Value = 1.436811
$ws.Range("O8").Value = 0.0813355572127976
$ws.Range("P8").Value = 0.08133555721279762
$ws.Range("Q8").Value = 0.01497891432066667
$ws.Range("R8").Value = 0.134810228886
$ws.Range("S8").Value = 0.01768658495789586
$ws.Range("T8").Value = 0.01768658495789587

# Row 9 (new)
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Fgf5"
$ws.Range("C9").Value = "Fgfr3"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.03127533333333334
$ws.Range("H9").Value = 0.09382600000000001
$ws.Range("I9").Value = 0.2174520660333412
$ws.Range("J9").Value = 0.2174520660333412
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.1116203333333333
$ws.Range("N9").Value = 0.334861
$ws.Range("O9").Value = 0.01895594202983873
$ws.Range("P9").Value = 0.01895594202983874
$ws.Range("Q9").Value = 0.003490963131777778
$ws.Range("R9").Value = 0.031418668186
$ws.Range("S9").Value = 0.004122008757996679
$ws.Range("T9").Value = 0.00412200875799668
